$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Suffrage) ---
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 0.6218629962506178
$ws.Range("D2").Value = 0.09729250143089825
$ws.Range("E2").Value = "qa_coverage_line_%"
$ws.Range("F2").Formula = "=1.851667786016466e-10"
$ws.Range("F2").Value2 = $ws.Range("F2").Value2

# --- Row 3 (GatesS) ---
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 0.4731027190065441
$ws.Range("D3").Value = 0.1231519890269868
$ws.Range("E3").Value = "qa_saccade_regression_rate_%"
$ws.Range("F3").Formula = "=7.472842657895719e-10"
$ws.Range("F3").Value2 = $ws.Range("F3").Value2

# --- Row 4 (GatesT) ---
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 0.7186502746103329
$ws.Range("D4").Value = 0.1263771545688306
$ws.Range("E4").Value = "qa_coverage_line_%"
$ws.Range("F4").Formula = "=6.794213183571248e-08"
$ws.Range("F4").Value2 = $ws.Range("F4").Value2
